# Auto-generated script applying scheduled market-data refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4585
$ws.Range("I64").Value = 4307.9
$ws.Range("J64").Value = 4892.8887
$ws.Range("K64").Value = 4307.9
$ws.Range("L64").Value = 4892.8887
$ws.Range("M64").Value = -4059.9
$ws.Range("N64").Value = -5388.8887
$ws.Range("H67").Value = 4585
$ws.Range("I67").Value = 4307.9
$ws.Range("J67").Value = 4892.8887
$ws.Range("K67").Value = 4307.9
$ws.Range("L67").Value = 4892.8887
$ws.Range("M67").Value = -3449.9
$ws.Range("N67").Value = -6608.8887
$ws.Range("H112").Value = 1574.8837
$ws.Range("J112").Value = 1675.5555
$ws.Range("L112").Value = 5026.666499999999
$ws.Range("N112").Value = -7242.666499999999
$ws.Range("H137").Value = 1720645.1
$ws.Range("I137").Value = 2022850.2
$ws.Range("J137").Value = 1364474.9
$ws.Range("K137").Value = 6068550.6
$ws.Range("L137").Value = 4093424.7
$ws.Range("M137").Value = -6066000.6
$ws.Range("N137").Value = -4098524.7
$ws.Range("H138").Value = 2949.1965
$ws.Range("I138").Value = 5088.8887
$ws.Range("J138").Value = 2539.468
$ws.Range("K138").Value = 15266.6661
$ws.Range("L138").Value = 7618.404
$ws.Range("M138").Value = -10126.6661
$ws.Range("N138").Value = -17898.404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 12483.875
$ws.Range("J31").Value = 23333.334
$ws.Range("L31").Value = 23333.334
$ws.Range("N31").Value = -23921.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 731.3
$ws.Range("I16").Value = 728.5714
$ws.Range("J16").Value = 737.6667
$ws.Range("K16").Value = 728.5714
$ws.Range("L16").Value = 737.6667
$ws.Range("M16").Value = -441.5714
$ws.Range("N16").Value = -1311.6667
$ws.Range("H31").Value = 1797.293
$ws.Range("I31").Value = 853.975
$ws.Range("J31").Value = 2436.8306
$ws.Range("K31").Value = 853.975
$ws.Range("L31").Value = 2436.8306
$ws.Range("M31").Value = -558.975
$ws.Range("N31").Value = -3026.8306
$ws.Range("H34").Value = 1797.293
$ws.Range("I34").Value = 853.975
$ws.Range("J34").Value = 2436.8306
$ws.Range("K34").Value = 853.975
$ws.Range("L34").Value = 2436.8306
$ws.Range("M34").Value = -651.975
$ws.Range("N34").Value = -2840.8306
$ws.Range("H99").Value = 85498.164
$ws.Range("I99").Value = 500506
$ws.Range("J99").Value = 2496.6
$ws.Range("K99").Value = 500506
$ws.Range("L99").Value = 2496.6
$ws.Range("M99").Value = -499008
$ws.Range("N99").Value = -5492.6
$ws.Range("H113").Value = 731.3
$ws.Range("I113").Value = 728.5714
$ws.Range("J113").Value = 737.6667
$ws.Range("K113").Value = 728.5714
$ws.Range("L113").Value = 737.6667
$ws.Range("M113").Value = 1441.4286
$ws.Range("N113").Value = -5077.6667
$ws.Range("H126").Value = 85498.164
$ws.Range("I126").Value = 500506
$ws.Range("J126").Value = 2496.6
$ws.Range("K126").Value = 1501518
$ws.Range("L126").Value = 7489.799999999999
$ws.Range("M126").Value = -1499048
$ws.Range("N126").Value = -12429.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1734
$ws.Range("I68").Value = 1304
$ws.Range("J68").Value = 1882.2759
$ws.Range("K68").Value = 3912
$ws.Range("L68").Value = 5646.8277
$ws.Range("M68").Value = -3101
$ws.Range("N68").Value = -7268.8277
$ws.Range("H71").Value = 1734
$ws.Range("I71").Value = 1304
$ws.Range("J71").Value = 1882.2759
$ws.Range("K71").Value = 11736
$ws.Range("L71").Value = 16940.4831
$ws.Range("M71").Value = -7680
$ws.Range("N71").Value = -25052.4831
$ws.Range("H107").Value = 1738.7084
$ws.Range("I107").Value = 836.6667
$ws.Range("J107").Value = 2279.9333
$ws.Range("K107").Value = 2510.0001
$ws.Range("L107").Value = 6839.7999
$ws.Range("M107").Value = -590.0001000000002
$ws.Range("N107").Value = -10679.7999
$ws.Range("H113").Value = 496.65714
$ws.Range("I113").Value = 541.65
$ws.Range("J113").Value = 436.66666
$ws.Range("K113").Value = 1624.95
$ws.Range("L113").Value = 1309.99998
$ws.Range("M113").Value = 545.0500000000002
$ws.Range("N113").Value = -5649.999980000001
$ws.Range("H120").Value = 2000
$ws.Range("I120").Value = 2000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 6000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -1162
$ws.Range("H121").Value = 3483.6428
$ws.Range("I121").Value = 3874.5557
$ws.Range("J121").Value = 2780
$ws.Range("K121").Value = 11623.6671
$ws.Range("L121").Value = 8340
$ws.Range("M121").Value = -10313.6671
$ws.Range("N121").Value = -10960
$ws.Range("H122").Value = 1008.1818
$ws.Range("I122").Value = 318.2
$ws.Range("J122").Value = 1308.174
$ws.Range("K122").Value = 2863.8
$ws.Range("L122").Value = 11773.566
$ws.Range("M122").Value = -413.7999999999997
$ws.Range("N122").Value = -16673.566
$ws.Range("H123").Value = 9486
$ws.Range("I123").Value = 30030
$ws.Range("J123").Value = 4350
$ws.Range("K123").Value = 90090
$ws.Range("L123").Value = 13050
$ws.Range("M123").Value = -87640
$ws.Range("N123").Value = -17950
$ws.Range("H124").Value = 1094.125
$ws.Range("I124").Value = 230
$ws.Range("J124").Value = 1958.25
$ws.Range("K124").Value = 690
$ws.Range("L124").Value = 5874.75
$ws.Range("M124").Value = 4220
$ws.Range("N124").Value = -15694.75
$ws.Range("H125").Value = 3950
$ws.Range("I125").Value = 3950
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 11850
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6930
$ws.Range("H126").Value = 7379.3
$ws.Range("I126").Value = 8837.143
$ws.Range("J126").Value = 3977.6667
$ws.Range("K126").Value = 26511.429
$ws.Range("L126").Value = 11933.0001
$ws.Range("M126").Value = -21571.429
$ws.Range("N126").Value = -21813.0001
$ws.Range("H127").Value = 5000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920
$ws.Range("H128").Value = 160000
$ws.Range("I128").Value = 160000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 480000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -475020
$ws.Range("H129").Value = 1354.3529
$ws.Range("I129").Value = 645.5
$ws.Range("J129").Value = 1850.55
$ws.Range("K129").Value = 1936.5
$ws.Range("L129").Value = 5551.65
$ws.Range("M129").Value = 3063.5
$ws.Range("N129").Value = -15551.65
$ws.Range("H130").Value = 112934.336
$ws.Range("I130").Value = 750
$ws.Range("J130").Value = 126957.375
$ws.Range("K130").Value = 2250
$ws.Range("L130").Value = 380872.125
$ws.Range("M130").Value = 2770
$ws.Range("N130").Value = -390912.125
$ws.Range("H131").Value = 1316.7548
$ws.Range("I131").Value = 2750.75
$ws.Range("J131").Value = 1199.6938
$ws.Range("K131").Value = 8252.25
$ws.Range("L131").Value = 3599.0814
$ws.Range("M131").Value = -3212.25
$ws.Range("N131").Value = -13679.0814
$ws.Range("H132").Value = 3888.889
$ws.Range("I132").Value = 4517.3335
$ws.Range("J132").Value = 2632
$ws.Range("K132").Value = 40656.0015
$ws.Range("L132").Value = 23688
$ws.Range("M132").Value = -38126.0015
$ws.Range("N132").Value = -28748
$ws.Range("H133").Value = 3239.3572
$ws.Range("I133").Value = 1255
$ws.Range("J133").Value = 4727.625
$ws.Range("K133").Value = 3765
$ws.Range("L133").Value = 14182.875
$ws.Range("M133").Value = 1295
$ws.Range("N133").Value = -24302.875
$ws.Range("H134").Value = 6603.625
$ws.Range("I134").Value = 4788.6
$ws.Range("J134").Value = 9628.666999999999
$ws.Range("K134").Value = 14365.8
$ws.Range("L134").Value = 28886.001
$ws.Range("M134").Value = -9295.800000000001
$ws.Range("N134").Value = -39026.001
$ws.Range("H136").Value = 2103.0625
$ws.Range("I136").Value = 1340.8182
$ws.Range("J136").Value = 3780
$ws.Range("K136").Value = 4022.4546
$ws.Range("L136").Value = 11340
$ws.Range("M136").Value = 1077.5454
$ws.Range("N136").Value = -21540
$ws.Range("H137").Value = 5733.4287
$ws.Range("I137").Value = 826.8
$ws.Range("J137").Value = 18000
$ws.Range("K137").Value = 2480.4
$ws.Range("L137").Value = 54000
$ws.Range("M137").Value = 2619.6
$ws.Range("N137").Value = -64200
$ws.Range("H138").Value = 1348.5714
$ws.Range("I138").Value = 1240
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 3720
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 1420
$ws.Range("N138").Value = -16280
$ws.Range("H139").Value = 2494.375
$ws.Range("I139").Value = 1109.1666
$ws.Range("J139").Value = 6650
$ws.Range("K139").Value = 3327.4998
$ws.Range("L139").Value = 19950
$ws.Range("M139").Value = 1812.5002
$ws.Range("N139").Value = -30230
$ws.Range("H140").Value = 1963.4375
$ws.Range("I140").Value = 1829.6428
$ws.Range("J140").Value = 2900
$ws.Range("K140").Value = 5488.928400000001
$ws.Range("L140").Value = 8700
$ws.Range("M140").Value = -308.9284000000007
$ws.Range("N140").Value = -19060
$ws.Range("H141").Value = 3768.024
$ws.Range("I141").Value = 1732.68
$ws.Range("J141").Value = 6761.1763
$ws.Range("K141").Value = 5198.04
$ws.Range("L141").Value = 20283.5289
$ws.Range("M141").Value = -18.03999999999996
$ws.Range("N141").Value = -30643.5289

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2740
$ws.Range("I7").Value = 2401
$ws.Range("J7").Value = 2966
$ws.Range("K7").Value = 2401
$ws.Range("L7").Value = 2966
$ws.Range("M7").Value = -2289
$ws.Range("N7").Value = -3190
$ws.Range("H8").Value = 2740
$ws.Range("I8").Value = 2401
$ws.Range("J8").Value = 2966
$ws.Range("K8").Value = 2401
$ws.Range("L8").Value = 2966
$ws.Range("M8").Value = -2262
$ws.Range("N8").Value = -3244
$ws.Range("H126").Value = 1867.4482
$ws.Range("I126").Value = 1818.4762
$ws.Range("J126").Value = 1996
$ws.Range("K126").Value = 5455.4286
$ws.Range("L126").Value = 5988
$ws.Range("M126").Value = -2985.4286
$ws.Range("N126").Value = -10928

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3151.3215
$ws.Range("I68").Value = 2654.8333
$ws.Range("J68").Value = 4045
$ws.Range("K68").Value = 2654.8333
$ws.Range("L68").Value = 4045
$ws.Range("M68").Value = -1905.8333
$ws.Range("N68").Value = -5543
$ws.Range("H71").Value = 3151.3215
$ws.Range("I71").Value = 2654.8333
$ws.Range("J71").Value = 4045
$ws.Range("K71").Value = 13274.1665
$ws.Range("L71").Value = 20225
$ws.Range("M71").Value = -9530.166499999999
$ws.Range("N71").Value = -27713
$ws.Range("H94").Value = 28500
$ws.Range("J94").Value = 28500
$ws.Range("L94").Value = 28500
$ws.Range("N94").Value = -29852
$ws.Range("H140").Value = 49607.25
$ws.Range("J140").Value = 49607.25
$ws.Range("L140").Value = 49607.25
$ws.Range("N140").Value = -59967.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 2666.6667
$ws.Range("J15").Value = 2666.6667
$ws.Range("L15").Value = 2666.6667
$ws.Range("N15").Value = -3242.6667
$ws.Range("H136").Value = 20402060
$ws.Range("I136").Value = 28085032
$ws.Range("J136").Value = 5844847.5
$ws.Range("K136").Value = 84255096
$ws.Range("L136").Value = 17534542.5
$ws.Range("M136").Value = -84252546
$ws.Range("N136").Value = -17539642.5
